# PDP - SDD Update
# Project Documentation Package Update: the "Unit & Integration Testing"
# slide title was previously split across two runs (an artifact of how it
# had been typed/edited in PowerPoint). Re-set the title text so it
# collapses back into a single run, the same way PowerPoint does when the
# text of a run is edited/retyped as a whole rather than appended to.

$p = $ppt.ActivePresentation

# Slide 8 is "Unit & Integration Testing" (Software Design Description /
# Unit & Integration Testing section of the Project Documentation Package).
$slide = $p.Slides.Item(8)
$title = $slide.Shapes.Title

$titleText = "Unit & Integration Testing"

# Using Replace (old text == new text) forces PowerPoint's text engine to
# re-flow/normalize the run structure for the paragraph, merging what were
# two runs ("Unit & " + "Integration Testing") back into a single run
# while preserving the existing run formatting - matching how the deck
# looks once the SDD section heading edit was finalized.
$title.TextFrame.TextRange.Replace($titleText, $titleText, 1, 0, 0) | Out-Null
